$wb = $excel.ActiveWorkbook
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$lastSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = '2025-09-09'

$newSheet.Cells.Item(2, 1).Value = 1
$newSheet.Cells.Item(2, 2).Value = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$newSheet.Cells.Item(2, 3).Value = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$newSheet.Cells.Item(2, 4).Value = '第５２話　暴走を止める器用貧乏（２）'
$newSheet.Cells.Item(3, 1).Value = 2
$newSheet.Cells.Item(3, 2).Value = '【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！'
$newSheet.Cells.Item(3, 3).Value = '島知宏 音速炒飯 有都あらゆる'
$newSheet.Cells.Item(3, 4).Value = '第２３食　巨大ヘビモンスターさん、パクパクですわ！（４）'
$newSheet.Cells.Item(4, 1).Value = 3
$newSheet.Cells.Item(4, 2).Value = '魔導具師ダリヤはうつむかない ～Dahliya Wilts No More～'
$newSheet.Cells.Item(4, 3).Value = '漫画：住川惠 原作：甘岸久弥(｢魔導具師ダリヤはうつむかない ～今日から自由な職人ライフ～｣MFブックス刊) キャラクター原案：景、駒田ハチ'
$newSheet.Cells.Item(4, 4).Value = '第47話 魔導具師とつながれたもの④'
$newSheet.Cells.Item(5, 1).Value = 4
$newSheet.Cells.Item(5, 2).Value = '十年目、帰還を諦めた転移者はいまさら主人公になる'
$newSheet.Cells.Item(5, 3).Value = '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう'
$newSheet.Cells.Item(5, 4).Value = '第１９話①'
$newSheet.Cells.Item(6, 1).Value = 5
$newSheet.Cells.Item(6, 2).Value = '塔の管理をしてみよう'
$newSheet.Cells.Item(6, 3).Value = '盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)'
$newSheet.Cells.Item(6, 4).Value = '第92話後編'
$newSheet.Cells.Item(7, 1).Value = 6
$newSheet.Cells.Item(7, 2).Value = '異世界でも無難に生きたい症候群'
$newSheet.Cells.Item(7, 3).Value = '原作：安泰（一二三書房刊） 漫画：笹峰コウ キャラクター原案：ひたきゆう'
$newSheet.Cells.Item(7, 4).Value = '第31話③'
$newSheet.Cells.Item(8, 1).Value = 7
$newSheet.Cells.Item(8, 2).Value = 'めっちゃ召喚された件 THE COMIC'
$newSheet.Cells.Item(8, 3).Value = '漫画：六甲島カモメ 原作：さいとうさ キャラクター原案：ツグトク'
$newSheet.Cells.Item(8, 4).Value = '第48話①'
$newSheet.Cells.Item(9, 1).Value = 8
$newSheet.Cells.Item(9, 2).Value = 'ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ'
$newSheet.Cells.Item(9, 3).Value = '漫画：晴野しゅー 原作：ちんくるり キャラクター原案：イセ川ヤスタカ'
$newSheet.Cells.Item(9, 4).Value = '第73話前半'
$newSheet.Cells.Item(10, 1).Value = 9
$newSheet.Cells.Item(10, 2).Value = 'レベル１だけどユニークスキルで最強です'
$newSheet.Cells.Item(10, 3).Value = '漫画：真綿 原作：三木なずな キャラクター原案：すばち'
$newSheet.Cells.Item(10, 4).Value = '第７４話　脱ブラックパーティー!? 転職のススメ！（１）'
$newSheet.Cells.Item(11, 1).Value = 10
$newSheet.Cells.Item(11, 2).Value = '王都ワンオペゴーレムマスター。まさかの追放！？～自由の身になったので弟子の美人勇者たちと一緒に最強ゴーレム作ります。戻ってこいと言われてももう知らん！～@COMIC'
$newSheet.Cells.Item(11, 3).Value = '阿住 周（漫画） レルクス（原作） 布施龍太（キャラクター原案）'
$newSheet.Cells.Item(11, 4).Value = '第10話'
$newSheet.Cells.Item(12, 1).Value = 11
$newSheet.Cells.Item(12, 2).Value = '骨ドラゴンのマナ娘'
$newSheet.Cells.Item(12, 3).Value = '雪白いち'
$newSheet.Cells.Item(12, 4).Value = '第39話「湯煙竜情②」'
$newSheet.Cells.Item(13, 1).Value = 12
$newSheet.Cells.Item(13, 2).Value = '外れスキル『レベルアップ』のせいでパーティーを追放された少年は、レベルを上げて物理で殴る'
$newSheet.Cells.Item(13, 3).Value = 'しんこせい 大橋ウルオ てんまそ'
$newSheet.Cells.Item(13, 4).Value = '第20話　パーティー（前編）'
$newSheet.Cells.Item(14, 1).Value = 13
$newSheet.Cells.Item(14, 2).Value = '人外姫様、始めました　-Free Life Fantasy Online-'
$newSheet.Cells.Item(14, 3).Value = '園原アオ 割田コマ 子日あきすず Sherry'
$newSheet.Cells.Item(14, 4).Value = '第６１話　冥府の王女？ それとも幽世の王女？（２）'
$newSheet.Cells.Item(15, 1).Value = 14
$newSheet.Cells.Item(15, 2).Value = 'ガヴリールドロップアウト'
$newSheet.Cells.Item(15, 3).Value = 'うかみ(著者)'
$newSheet.Cells.Item(15, 4).Value = '第127話'
$newSheet.Cells.Item(16, 1).Value = 15
$newSheet.Cells.Item(16, 2).Value = '無能と呼ばれた『精霊たらし』～実は異能で、精霊界では伝説的ヒーローでした～＠COMIC'
$newSheet.Cells.Item(16, 3).Value = '原作：佐藤謙羊 漫画：タバタグランドキャニオン'
$newSheet.Cells.Item(16, 4).Value = '第29話「カレキット村の奇跡」②'
$newSheet.Cells.Item(17, 1).Value = 16
$newSheet.Cells.Item(17, 2).Value = '器用貧乏、城を建てる～開拓学園の劣等生なのに、上級職のスキルと魔法がすべて使えます～＠COMIC'
$newSheet.Cells.Item(17, 3).Value = '原作：佐藤謙羊 漫画：スガン'
$newSheet.Cells.Item(17, 4).Value = '第23話③「自爆スイッチは押されたい」'
$newSheet.Cells.Item(18, 1).Value = 17
$newSheet.Cells.Item(18, 2).Value = '国王である兄から辺境に追放されたけど平穏に暮らしたい ～目指せスローライフ～'
$newSheet.Cells.Item(18, 3).Value = 'おとら(原作) 西沢秀二(漫画) 夜ノみつき(キャラクター原案)'
$newSheet.Cells.Item(18, 4).Value = '第11話-2'
$newSheet.Cells.Item(19, 1).Value = 18
$newSheet.Cells.Item(19, 2).Value = 'ブチ切れ令嬢は報復を誓いました。 ～魔導書の力で祖国を叩き潰します～'
$newSheet.Cells.Item(19, 3).Value = '漫画：おおのいも 原作：はぐれメタボ キャラクター原案：昌未'
$newSheet.Cells.Item(19, 4).Value = '第51話'
$newSheet.Cells.Item(20, 1).Value = 19
$newSheet.Cells.Item(20, 2).Value = '暴食のベルセルク～俺だけレベルという概念を突破する～'
$newSheet.Cells.Item(20, 3).Value = '漫画：滝乃大祐 原作：一色一凛 キャラクター原案：fame'
$newSheet.Cells.Item(20, 4).Value = '第74話前半'
$newSheet.Cells.Item(21, 1).Value = 20
$newSheet.Cells.Item(21, 2).Value = '転生少女はまず一歩からはじめたい～魔物がいるとか聞いてない！～'
$newSheet.Cells.Item(21, 3).Value = '原作：カヤ 漫画：岡村アユム キャラクター原案：那流'
$newSheet.Cells.Item(21, 4).Value = '第39歩 薬師修行はじまります①'
$newSheet.Cells.Item(22, 1).Value = 21
$newSheet.Cells.Item(22, 2).Value = 'まんきつしたい常連さん'
$newSheet.Cells.Item(22, 3).Value = 'しんみりん(著者)'
$newSheet.Cells.Item(22, 4).Value = '第47話後編'
$newSheet.Cells.Item(23, 1).Value = 22
$newSheet.Cells.Item(23, 2).Value = 'クラス最安値で売られた俺は、実は最強パラメーター'
$newSheet.Cells.Item(23, 3).Value = 'カンブリア爆発太郎(漫画) RYOMA(原作) 黒井ススム(キャラクター原案)'
$newSheet.Cells.Item(23, 4).Value = '第37話-2'
$newSheet.Cells.Item(24, 1).Value = 23
$newSheet.Cells.Item(24, 2).Value = '魔眼の悪役に転生したので推しキャラを見守るモブを目指します'
$newSheet.Cells.Item(24, 3).Value = '在間りしん(漫画) 瀧岡くるじ(原作) 福きつね(キャラクター原案)'
$newSheet.Cells.Item(24, 4).Value = '第12話②'
$newSheet.Cells.Item(25, 1).Value = 24
$newSheet.Cells.Item(25, 2).Value = '「門番やってろ」と言われ15年、突っ立ってる間に俺の魔力が9999（最強）に育ってました'
$newSheet.Cells.Item(25, 3).Value = '漫画：はり太郎 原作：まさキチ キャラクター原案：カラスBTK'
$newSheet.Cells.Item(25, 4).Value = '第3話'
$newSheet.Cells.Item(26, 1).Value = 25
$newSheet.Cells.Item(26, 2).Value = '真の聖女である私は追放されました。だからこの国はもう終わりです'
$newSheet.Cells.Item(26, 3).Value = '松もくば 鬱沢色素 ぷきゅのすけ'
$newSheet.Cells.Item(26, 4).Value = '第52話　女神との交信です！（２）'
$newSheet.Cells.Item(27, 1).Value = 26
$newSheet.Cells.Item(27, 2).Value = '不遇皇子は天才錬金術師～皇帝なんて柄じゃないので弟妹を可愛がりたい～@COMIC'
$newSheet.Cells.Item(27, 3).Value = '長先ザワ（漫画） うめー（原作） 瑛来イチ（構成） 雨銛（構成） かわく（キャラクター原案）'
$newSheet.Cells.Item(27, 4).Value = '第9話 ②'
$newSheet.Cells.Item(28, 1).Value = 27
$newSheet.Cells.Item(28, 2).Value = '願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜'
$newSheet.Cells.Item(28, 3).Value = 'ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)'
$newSheet.Cells.Item(28, 4).Value = '第6話-1：火蓮の剣'
$newSheet.Cells.Item(29, 1).Value = 28
$newSheet.Cells.Item(29, 2).Value = '狂戦士なモブ、無自覚に本編を破壊する'
$newSheet.Cells.Item(29, 3).Value = '漫画：佐藤良亮 原作：なるのるな キャラクター原案：霜月えいと'
$newSheet.Cells.Item(29, 4).Value = '第12話 ③'
$newSheet.Cells.Item(30, 1).Value = 29
$newSheet.Cells.Item(30, 2).Value = '２度目の人生、と思ったら、実は３度目だった。～歴史知識と内政努力で不幸な歴史の改変に挑みます～@COMIC'
$newSheet.Cells.Item(30, 3).Value = '麦こうちゃ（漫画） take4（原作） 桧野ひなこ（キャラクター原案）'
$newSheet.Cells.Item(30, 4).Value = '第9話 ②'
$newSheet.Cells.Item(31, 1).Value = 30
$newSheet.Cells.Item(31, 2).Value = 'マジカル★エクスプローラー エロゲの友人キャラに転生したけど、ゲーム知識使って自由に生きる'
$newSheet.Cells.Item(31, 3).Value = '入栖(原作) 緋賀ゆかり(漫画) 神奈月 昇(キャラクター原案)'
$newSheet.Cells.Item(31, 4).Value = '第3話-1'
$newSheet.Cells.Item(32, 1).Value = 31
$newSheet.Cells.Item(32, 2).Value = '異世界のんびり開拓記  -平凡サラリーマン、万能自在のビルド&クラフトスキルで、気ままなスローライフ 開拓始めます! -'
$newSheet.Cells.Item(32, 3).Value = '漫画：しょうじひでまさ 原作：タライ和治 キャラクター原案：イシバシヨウスケ'
$newSheet.Cells.Item(32, 4).Value = '第24話'
$newSheet.Cells.Item(33, 1).Value = 32
$newSheet.Cells.Item(33, 2).Value = 'SSSランクダンジョンでナイフ一本手渡され追放された白魔導師 ユグドラシルの呪いにより弱点である魔力不足を克服し世界最強へと至る'
$newSheet.Cells.Item(33, 3).Value = '上下瑞樹(漫画) カミトイチ(原作) 眠介(キャラクター原案)'
$newSheet.Cells.Item(33, 4).Value = '第23話-3'
$newSheet.Cells.Item(34, 1).Value = 33
$newSheet.Cells.Item(34, 2).Value = 'S級パーティーから追放された狩人、実は世界最強 ～射程9999の男、帝国の狙撃手として無双する～'
$newSheet.Cells.Item(34, 3).Value = '漫画：カズミヤアキラ 原作：茨木野 キャラクター原案：へいろー'
$newSheet.Cells.Item(34, 4).Value = '第10話 ③'
$newSheet.Cells.Item(35, 1).Value = 34
$newSheet.Cells.Item(35, 2).Value = 'コボルト無双、モフモフな最弱噛ませ犬だけど世界最強を目指す！'
$newSheet.Cells.Item(35, 3).Value = '赤志木ひの乃 shiba'
$newSheet.Cells.Item(35, 4).Value = '第十六話 夜盗との戦い'
$newSheet.Cells.Item(36, 1).Value = 35
$newSheet.Cells.Item(36, 2).Value = 'クセ強彼女は床にいざなう'
$newSheet.Cells.Item(36, 3).Value = '須河篤志(著者)'
$newSheet.Cells.Item(36, 4).Value = '第15話前半'
$newSheet.Cells.Item(37, 1).Value = 36
$newSheet.Cells.Item(37, 2).Value = '能あるオーガは角を隠す'
$newSheet.Cells.Item(37, 3).Value = '漫画家： 蒼葉 結 原作： 津野瀬 文'
$newSheet.Cells.Item(37, 4).Value = '第10話 後編'
$newSheet.Cells.Item(38, 1).Value = 37
$newSheet.Cells.Item(38, 2).Value = '無能は不要と言われ『時計使い』の僕は職人ギルドから追い出されるも、ダンジョンの深部で真の力に覚醒する'
$newSheet.Cells.Item(38, 3).Value = '漫画：さらさみさ 小説： 桜霧琥珀 キャラ原案： 福きつね'
$newSheet.Cells.Item(38, 4).Value = '第20話前半'
$newSheet.Cells.Item(39, 1).Value = 38
$newSheet.Cells.Item(39, 2).Value = '追放されたギルド職員は、世界最強の召喚士@COMIC'
$newSheet.Cells.Item(39, 3).Value = '原作：月島秀一 漫画：あづち涼 キャラクター原案：チワワ丸'
$newSheet.Cells.Item(39, 4).Value = '第11話④「絶望の復魔十使」'
$newSheet.Cells.Item(40, 1).Value = 39
$newSheet.Cells.Item(40, 2).Value = 'ちはるくんは女装をしたくない！'
$newSheet.Cells.Item(40, 3).Value = '翁丸ジョン'
$newSheet.Cells.Item(40, 4).Value = '【第23話】男装宗と交流したくない！その五'
$newSheet.Cells.Item(41, 1).Value = 40
$newSheet.Cells.Item(41, 2).Value = 'ぽんドロイド！ はまさん'
$newSheet.Cells.Item(41, 3).Value = 'はれやまはれぞう(著者)'
$newSheet.Cells.Item(41, 4).Value = '第8話'
$newSheet.Cells.Item(42, 1).Value = 41
$newSheet.Cells.Item(42, 2).Value = '異世界転移で女神様から祝福を！～いえ、手持ちの異能があるので結構です～@COMIC'
$newSheet.Cells.Item(42, 3).Value = 'コーダ 壁アキオ'
$newSheet.Cells.Item(42, 4).Value = '第3話「盗賊の発見と退治①」'
$newSheet.Cells.Item(43, 1).Value = 42
$newSheet.Cells.Item(43, 2).Value = 'アラフォーおっさんはスローライフの夢を見るか？'
$newSheet.Cells.Item(43, 3).Value = '漫画：大関詠詞 原作：サイトウアユム キャラクター原案： ジョンディー'
$newSheet.Cells.Item(43, 4).Value = '第17話'
$newSheet.Cells.Item(44, 1).Value = 43
$newSheet.Cells.Item(44, 2).Value = 'まったく最近の探偵ときたら'
$newSheet.Cells.Item(44, 3).Value = '五十嵐正邦(著者)'
$newSheet.Cells.Item(44, 4).Value = '第115話'
$newSheet.Cells.Item(45, 1).Value = 44
$newSheet.Cells.Item(45, 2).Value = '氷結令嬢さまをフォローしたら、メチャメチャ溺愛されてしまった件@comic'
$newSheet.Cells.Item(45, 3).Value = '漫画：ハレノチアメ 原作：愛坂タカト キャラクター原案：Bcoca'
$newSheet.Cells.Item(45, 4).Value = 'アリシア様セクシーショット集（担当編集厳選）'
$newSheet.Cells.Item(46, 1).Value = 45
$newSheet.Cells.Item(46, 2).Value = '灰原くんの強くて青春ニューゲーム ヨコ読み版'
$newSheet.Cells.Item(46, 3).Value = '漫画：みさおまる、プラス81 原作：雨宮和希 キャラクター原案：吟'
$newSheet.Cells.Item(46, 4).Value = '第5話'
$newSheet.Cells.Item(47, 1).Value = 46
$newSheet.Cells.Item(47, 2).Value = 'リビルドワールド'
$newSheet.Cells.Item(47, 3).Value = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$newSheet.Cells.Item(47, 4).Value = '第72話④'
$newSheet.Cells.Item(48, 1).Value = 47
$newSheet.Cells.Item(48, 2).Value = '千年英雄'
$newSheet.Cells.Item(48, 3).Value = '原作/福島航平 作画/中村ゆきひろ'
$newSheet.Cells.Item(48, 4).Value = '22話②'
$newSheet.Cells.Item(49, 1).Value = 48
$newSheet.Cells.Item(49, 2).Value = 'ワンパンマン'
$newSheet.Cells.Item(49, 3).Value = '原作/ＯＮＥ 作画/村田雄介'
$newSheet.Cells.Item(49, 4).Value = '210撃目'
$newSheet.Cells.Item(50, 1).Value = 49
$newSheet.Cells.Item(50, 2).Value = '王子様の友達'
$newSheet.Cells.Item(50, 3).Value = 'すけろく(著者)'
$newSheet.Cells.Item(50, 4).Value = '第30話'
$newSheet.Cells.Item(51, 1).Value = 50
$newSheet.Cells.Item(51, 2).Value = 'わたしのために脱ぎなさいっ！'
$newSheet.Cells.Item(51, 3).Value = '九郎(著者)'
$newSheet.Cells.Item(51, 4).Value = '第85話'
